# Update the "signal segment 8" / "signal segment 9" rows (rows 9-10) across
# the Step1_Data / Step2_Sj / Step3_DataPts_* sheets to reflect the recomputed
# "mounted pipeline" values (the raw D9 bin mass moved into the new AJ9 bin,
# the rest of the Step1_Data row was renormalized, Step2_Sj's cumulative
# column and the downstream Step3 threshold lookups were refreshed to match).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.1500041412562959
$ws.Range("F9").Value = 0.1033480244528977
$ws.Range("G9").Value = 0.03612947664136002
$ws.Range("H9").Value = 0.03997766674740796
$ws.Range("I9").Value = 0.03850976266946552
$ws.Range("J9").Value = 0.003103880820280215
$ws.Range("K9").Value = 0.008132972574936037
$ws.Range("L9").Value = 0.1454709062747589
$ws.Range("M9").Value = 0.03741832765358959
$ws.Range("N9").Value = 0.01224415514374681
$ws.Range("O9").Value = 0.1247361464729475
$ws.Range("P9").Value = 0.004852600868462021
$ws.Range("T9").Value = 0.03692634152346212
$ws.Range("U9").Value = 0.06043404877168231
$ws.Range("V9").Value = 0.02532473082869448
$ws.Range("Y9").Value = 0.003014820928716276
$ws.Range("AA9").Value = 0.02077852908042098
$ws.Range("AC9").Value = 0.005818055617109732
$ws.Range("AD9").Value = 0.05954045807025455
$ws.Range("AE9").Value = 0.06571209671530268
$ws.Range("AF9").Value = 0.01119020730041218
$ws.Range("AI9").Value = 0.004616482950805654
$ws.Range("AJ9").Value = 0.002716166636991094
$ws.Range("E10").Value = 0.1573657228656474
$ws.Range("F10").Value = 0.1057821027733479
$ws.Range("I10").Value = 0.04600923688284865
$ws.Range("M10").Value = 0.05614301285459129
$ws.Range("S10").Value = 0.00309929032481564
$ws.Range("T10").Value = 0.02805691009604994
$ws.Range("U10").Value = 0.05215006988203307
$ws.Range("V10").Value = 0.03728132718516838
$ws.Range("AD10").Value = 0.05875949257041024
$ws.Range("AE10").Value = 0.06492364734926986

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0.1500041412562959
$ws.Range("F9").Value = 0.2533521657091936
$ws.Range("G9").Value = 0.2894816423505536
$ws.Range("H9").Value = 0.3294593090979616
$ws.Range("I9").Value = 0.3679690717674272
$ws.Range("J9").Value = 0.3710729525877074
$ws.Range("K9").Value = 0.3792059251626435
$ws.Range("L9").Value = 0.5246768314374024
$ws.Range("M9").Value = 0.5620951590909919
$ws.Range("N9").Value = 0.5743393142347387
$ws.Range("O9").Value = 0.6990754607076862
$ws.Range("P9").Value = 0.7039280615761483
$ws.Range("Q9").Value = 0.7039280615761483
$ws.Range("R9").Value = 0.7039280615761483
$ws.Range("S9").Value = 0.7039280615761483
$ws.Range("T9").Value = 0.7408544030996104
$ws.Range("U9").Value = 0.8012884518712927
$ws.Range("V9").Value = 0.8266131826999872
$ws.Range("W9").Value = 0.8266131826999872
$ws.Range("X9").Value = 0.8266131826999872
$ws.Range("Y9").Value = 0.8296280036287035
$ws.Range("Z9").Value = 0.8296280036287035
$ws.Range("AA9").Value = 0.8504065327091245
$ws.Range("AB9").Value = 0.8504065327091245
$ws.Range("AC9").Value = 0.8562245883262342
$ws.Range("AD9").Value = 0.9157650463964887
$ws.Range("AE9").Value = 0.9814771431117914
$ws.Range("AF9").Value = 0.9926673504122036
$ws.Range("AG9").Value = 0.9926673504122036
$ws.Range("AH9").Value = 0.9926673504122036
$ws.Range("AI9").Value = 0.9972838333630092
$ws.Range("E10").Value = 0.1573657228656474
$ws.Range("G10").Value = 0.281651624640451
$ws.Range("H10").Value = 0.3091668341064557
$ws.Range("I10").Value = 0.3551760709893044
$ws.Range("J10").Value = 0.3551760709893044
$ws.Range("K10").Value = 0.3551760709893044
$ws.Range("L10").Value = 0.5273969665266303
$ws.Range("M10").Value = 0.5835399793812216
$ws.Range("N10").Value = 0.5865653099581652
$ws.Range("O10").Value = 0.7420776278967631
$ws.Range("P10").Value = 0.7420776278967631
$ws.Range("Q10").Value = 0.7420776278967631
$ws.Range("R10").Value = 0.7420776278967631
$ws.Range("S10").Value = 0.7451769182215787
$ws.Range("T10").Value = 0.7732338283176287
$ws.Range("U10").Value = 0.8253838981996617
$ws.Range("V10").Value = 0.8626652253848301
$ws.Range("W10").Value = 0.8626652253848301
$ws.Range("X10").Value = 0.8626652253848301
$ws.Range("Y10").Value = 0.8626652253848301
$ws.Range("Z10").Value = 0.8626652253848301
$ws.Range("AA10").Value = 0.863905473584168
$ws.Range("AB10").Value = 0.863905473584168
$ws.Range("AC10").Value = 0.863905473584168
$ws.Range("AD10").Value = 0.9226649661545783
$ws.Range("AE10").Value = 0.9875886135038481
$ws.Range("AF10").Value = 0.9999999999999999
$ws.Range("AG10").Value = 0.9999999999999999
$ws.Range("AH10").Value = 0.9999999999999999
$ws.Range("AI10").Value = 0.9999999999999999
$ws.Range("AJ10").Value = 0.9999999999999999
$ws.Range("AK10").Value = 0.9999999999999999

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F9").Value = 0.5246768314374024
$ws.Range("F10").Value = 0.5273969665266303

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("D9").Value = 15
$ws.Range("F9").Value = 0.7039280615761483
$ws.Range("G9").Value = 13
$ws.Range("F10").Value = 0.7420776278967631

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F9").Value = 0.8012884518712927
$ws.Range("F10").Value = 0.8253838981996617

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F9").Value = 0.9157650463964887
$ws.Range("F10").Value = 0.9226649661545783
